# Resort the worksheets: move "总计" (the summary sheet) in front of
# "2020-Q4" (the detail sheet), so that "总计" becomes the first / active
# sheet and "2020-Q4" becomes the second sheet - swapping their order.

$wb = $excel.ActiveWorkbook

$detailSheet  = $wb.Worksheets.Item("2020-Q4")
$summarySheet = $wb.Worksheets.Item("总计")

# Move the summary sheet so that it sits right before the detail sheet.
$summarySheet.Move($detailSheet)
